$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data content of row 2 and row 3 for the columns that
# differ between the two rows (A, B, D, E, F, G, H, K, Q, R, AH, AJ, AK, AO).
# Columns that already hold identical values in both rows are left untouched.

$cols = @("A","B","D","E","F","G","H","K","Q","R","AH","AJ","AK","AO")

foreach ($col in $cols) {
    $addr2 = "$col`2"
    $addr3 = "$col`3"
    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}
